$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new data row 11 (mirrors rows 6-10 pattern)
$ws.Range("C11").Formula = "=1-(0.5)^(F11/E11)"
$ws.Range("D11").Formula = "=F11*(1-0.5^(1/E11))"
$ws.Range("E11").Value = 360
$ws.Range("F11").Value = 20

# Update the selected cell to match the new location (D12)
$ws.Range("D12").Select()
